$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Developed forecasting / updated data: append three new daily rows
# (95, 96, 97) below the existing combat-losses table, continuing
# the date series and per-column running totals.
# ---------------------------------------------------------------

# 1) Copy the formatting (styles) of the last existing row (94) down
#    onto the three new rows so the new cells keep the same look
#    (date format on column A, numeric/centered format on B:O).
$ws.Range("A94:O94").Copy()
$ws.Range("A95:O97").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 2) Restore the custom row height used throughout the table.
$ws.Rows.Item(95).RowHeight = 15.75
$ws.Rows.Item(96).RowHeight = 15.75
$ws.Rows.Item(97).RowHeight = 15.75

# 3) Column A keeps the running "date+1" formula.
$ws.Range("A95").Formula = "=A94+1"
$ws.Range("A96").Formula = "=A95+1"
$ws.Range("A97").Formula = "=A96+1"

# 4) Fill in the new data (columns B:O) for the three new rows.
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O")

$row95 = @(30000, 1330, 3258, 628, 203, 93, 207, 174, 2226, 13, 503, 48, 4, 116)
$row96 = @(30150, 1338, 3270, 631, 203, 93, 207, 174, 2240, 13, 504, 48, 4, 116)
$row97 = @(30350, 1349, 3282, 643, 205, 93, 207, 174, 2258, 13, 507, 48, 4, 118)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "95").Value2 = $row95[$i]
    $ws.Range($cols[$i] + "96").Value2 = $row96[$i]
    $ws.Range($cols[$i] + "97").Value2 = $row97[$i]
}

# 5) Update the view: scroll near the bottom of the table and leave
#    the selection on the newly added data (mirrors the author's
#    last cursor position after entering the new rows).
[void]$excel.Goto($ws.Range("M80"), $true)
[void]$ws.Range("M106").Select()
